# Update the "Förändrad" (Changed) date column (C) from 46074 to 46075
# for every data row (rows 2 through 36) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 46074) {
        $cell.Value = 46075
    }
}
